$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workload")

# Row 5: student IDs (numeric, previously "<student id>" placeholder text)
$ws.Range("D5").Value = 5762340
$ws.Range("E5").Value = 5556910
$ws.Range("F5").Value = 5699193

# Row 6: student names (previously "<student name>" placeholder text).
# Write F/E/D in this order so the new shared strings are appended in the
# same order as the target workbook (Andrea, Eduard, Tudor).
$ws.Range("F6").Value = "Andrea Vezzuto"
$ws.Range("E6").Value = "Eduard Faraon"
$ws.Range("D6").Value = "Tudor Coman"

# Basic features section (rows 8-15): flip several 0% entries to 100%
$ws.Range("D8").Value = 100
$ws.Range("E9").Value = 100
$ws.Range("E10").Value = 100
$ws.Range("F11").Value = 100
$ws.Range("E12").Value = 100
$ws.Range("E13").Value = 100
$ws.Range("F14").Value = 100
$ws.Range("E15").Value = 100

# Extra features section (rows 19-24): flip several 0% entries to 100%
$ws.Range("D19").Value = 100
$ws.Range("D20").Value = 100
$ws.Range("F21").Value = 100
$ws.Range("F22").Value = 100
$ws.Range("F23").Value = 100
$ws.Range("E24").Value = 100

# Restore the selection to match where the author left the cursor on save.
[void]$ws.Range("J37").Select()
